$p = $ppt.ActivePresentation
Write-Output ($p | Get-Member | Out-String)
